# Update crypto price/volume data per the Sun Apr 28 20:29:38 UTC 2024 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text format
# so Excel keeps them as strings (preserving formatting such as trailing zeros)
# instead of silently converting them to numeric values.
$textCells = @("D5", "D6", "D11", "D14", "D19", "D20", "D24", "D25", "D27", "D29", "D32", "D35", "D36", "D37", "D39", "D41", "D43", "D46", "D47", "D48", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.683.18'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '3.313.36'
$ws.Range("E3").Value = '  +2.66%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '606.61'
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").Value = '141.55'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.312.12'
$ws.Range("E8").Value = '  +2.70%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("D11").Value = '5.53'
$ws.Range("E11").Value = '  +3.37%  '
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").Value = '34.95'
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").Value = '3.858.88'
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '3.313.77'
$ws.Range("E17").Value = '  +2.54%  '
$ws.Range("D18").Value = '63.761.70'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("D20").Value = '480.70'
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("E22").Value = '  +1.98%  '
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("D24").Value = '14.00'
$ws.Range("E24").Value = '  +6.64%  '
$ws.Range("D25").Value = '85.17'
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  +1.48%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '8.22'
$ws.Range("E29").Value = '  +1.66%  '
$ws.Range("E30").Value = '  -4.48%  '
$ws.Range("E31").Value = '  +2.22%  '
$ws.Range("D32").Value = '28.87'
$ws.Range("E32").Value = '  +5.56%  '
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  +1.97%  '
$ws.Range("D36").Value = '6.06'
$ws.Range("E36").Value = '  +2.70%  '
$ws.Range("D37").Value = '52.50'
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").Value = '0.0₃0748'
$ws.Range("E38").Value = '  +5.88%  '
$ws.Range("D39").Value = '0.0400'
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("D40").Value = '3.119.05'
$ws.Range("E40").Value = '  +4.94%  '
$ws.Range("D41").Value = '432.96'
$ws.Range("E41").Value = '  +2.54%  '
$ws.Range("E42").Value = '  +9.04%  '
$ws.Range("D43").Value = '8.35'
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = '2.23'
$ws.Range("E46").Value = '  +3.31%  '
$ws.Range("D47").Value = '36.82'
$ws.Range("E47").Value = '  +8.84%  '
$ws.Range("D48").Value = '26.29'
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("D50").Value = '2.32'
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("E51").Value = '  -0.48%  '
